# Applies the LegacyCollectionsTemplate content edit:
#  - "Captain: {{" / "captain-name" / "}}" (3 runs) -> single run
#    "Captain: {{captain-name}}"
#  - new paragraph inserted right after the Captain paragraph containing
#    "{{" + (bookmarked) "c" + "aptain-first-officer-name}}"
#  - table placeholders for officer-name / officer-rank / officer-uniform
#    collapsed from 3 runs ("{{" / "name" / "}}") into a single run each
#  - "{{" / "service-ship-name" / "}}" (3 runs) -> single run
#    "{{service-ship-name}}"

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1 (unused here), wdReplaceAll = 2
# wdCollapseEnd = 0, wdCharacter = 1

# --- 1. Merge the "Captain: {{" + "captain-name" + "}}" runs ------------
$find = $d.Content.Find
$find.Execute("Captain: {{captain-name}}", $false, $false, $false, $false, `
               $false, $true, 1, $false, "Captain: {{captain-name}}", 2)

# --- 2. Insert a new paragraph right after it holding the new ------------
#        captain-first-officer-name placeholder, wrapping a bookmark
#        around the single character "c" (matching the authored diff).
$captainPara = $d.Paragraphs.Item(1)
$afterRange = $captainPara.Range
$afterRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item(2)
$newRange = $newPara.Range
$newRange.MoveEnd(1, -1)
$newRange.Text = "{{captain-first-officer-name}}"

$bmStart = $newPara.Range.Start + 2
$bmRange = $d.Range($bmStart, $bmStart + 1)
$d.Bookmarks.Add("__DdeLink__44_427613671", $bmRange)

# --- 3. Table placeholders: officer-name / officer-rank / officer-uniform
$officerTable = $d.Tables.Item(1)

$cellName = $officerTable.Cell(2, 1)
$cellNameFind = $cellName.Range.Find
$cellNameFind.Execute("{{officer-name}}", $false, $false, $false, $false, `
                       $false, $true, 1, $false, "{{officer-name}}", 2)

$cellRank = $officerTable.Cell(2, 2)
$cellRankFind = $cellRank.Range.Find
$cellRankFind.Execute("{{officer-rank}}", $false, $false, $false, $false, `
                       $false, $true, 1, $false, "{{officer-rank}}", 2)

$cellUniform = $officerTable.Cell(2, 3)
$cellUniformFind = $cellUniform.Range.Find
$cellUniformFind.Execute("{{officer-uniform}}", $false, $false, $false, $false, `
                          $false, $true, 1, $false, "{{officer-uniform}}", 2)

# --- 4. Merge "{{" + "service-ship-name" + "}}" runs ---------------------
$find2 = $d.Content.Find
$find2.Execute("{{service-ship-name}}", $false, $false, $false, $false, `
                $false, $true, 1, $false, "{{service-ship-name}}", 2)
